# 7.8 History Card & Advanced Story
# Updates the Lee interrogation dialogue rows: fixes the "Steward"/"Butler"
# naming, doubles up the em-dash interruptions, and wraps the sabers clue
# line in a green color tag. Shared strings that are no longer referenced
# by any cell are compacted away by the engine (same as native Excel would
# do over repeated edits), and brand-new strings are appended to the shared
# string table in the order they are first assigned below - so the order
# of assignment below is chosen to reproduce the target shared-string
# table layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose text is unchanged (just keep shared strings alive) ---
$ws.Range("B4").Value  = "Don’t worry, I know you’re innocent."
$ws.Range("B7").Value  = "When was the last time you saw the Lord?"
$ws.Range("B8").Value  = "It should’ve been after lunch. I ran into the Lord in the corridor."
$ws.Range("B9").Value  = "He kindly asked about my injury and told me not to leave the manor until I’d fully recovered."
$ws.Range("B10").Value = "I didn’t see him again after that."
$ws.Range("B11").Value = "Tell me about what you did before and during the banquet."
$ws.Range("B13").Value = "My leg hadn’t healed yet, and Doctor Ran had just changed my bandages, so I wanted to arrive early and find a seat."
$ws.Range("B14").Value = "Was anyone else there when you arrived at the banquet hall?"
$ws.Range("B15").Value = "No, I was the first one there."
$ws.Range("B16").Value = "But Lai arrived shortly after, and then you and your attendant came."
$ws.Range("B17").Value = "Oh right, I also have a clue I want to share."
$ws.Range("B18").Value = "What is it?"
$ws.Range("B19").Value = "I presented the twin sabers I forged to the Lord this afternoon, and he immediately strapped them on."
$ws.Range("B20").Value = "But when I secretly examined the body just now, I found the sabers were missing."
$ws.Range("B21").Value = "Alright. Thank you."

# --- Row 22: wrap the existing clue line in a green color tag ---
$ws.Range("B22").Value = " <color=#00CC00>(It’s unlikely that stealing weapons was the killer’s motive, but the missing sabers are certainly an important clue.)</color>"

# --- Rows whose text actually changes (doubled em-dash / Steward->Butler) ---
$ws.Range("B3").Value  = "Huh? Don’t suspect me——I injured my leg! There’s no way I could be the killer!"
$ws.Range("B5").Value  = "You arrived at the banquet hall early and didn’t leave at all during the dinner——I can vouch for that."
$ws.Range("B6").Value  = "That’s a relief. Ask anything you like——I’ll answer truthfully."
$ws.Range("B12").Value = "After Butler He came to inform me of the time and place, I headed out immediately."

# --- View state: reflect the post-edit selection ---
$ws.Range("B24").Select()
